# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row at 33 so everything below (old rows 33-38) shifts
#    down by one (old row37 -> row38, old row38 -> row39). The row that
#    used to be the last detail line (old row32, period 2403) will stay
#    at row32 but be re-purposed into a normal detail row, and the new
#    blank row33 becomes the new last detail line (period 2508).
# ---------------------------------------------------------------------
$ws.Rows(33).Insert()

# ---------------------------------------------------------------------
# 2) Re-apply correct formatting:
#    - row33 should look like the old "last row" (heavier bottom border)
#      that row32 still has right now, so copy row32's format down.
#    - row32 should become a normal detail row like rows 16-31.
# ---------------------------------------------------------------------
$ws.Range("B32:J32").Copy()
$ws.Range("B33:J33").PasteSpecial(-4122)

$ws.Range("B31:J31").Copy()
$ws.Range("B32:J32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Fill in the detail values. The account now lists 18 consecutive
#    periods (2403 .. 2508, ascending) instead of the previous 17
#    (2403 .. 2507 descending in rows 16-32). Row16 is the new oldest
#    period (partial month -> lower value), row33 is the brand-new
#    last period (2508).
# ---------------------------------------------------------------------
$periods = @("2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "45741396"
    $ws.Range("D$r").Value = "ELSY FUENTES NIETO"
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("G$r").Value = 650000
}
$ws.Range("F16").Value = 8667
for ($r = 17; $r -le 33; $r++) {
    $ws.Range("F$r").Value = 26000
}

# ---------------------------------------------------------------------
# 4) Header / summary updates.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 450667
$ws.Range("F13").Value = 18
